# SOLUSDT_analysis.xlsx update: append 5 new scrape rows (23-27 on "Data";
# 25-29 on the per-topic sheets), refresh the Dashboard "latest analysis"
# summary cells, and extend every chart series range from row 24 to row 29.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New rows of scraped data (Data sheet rows 23-27), columns A..AN.
# ---------------------------------------------------------------------
$dataRows = @(
    @('16:40:09 21/01/2026', 46043.6945559778, 127.56, 129.77, 124.68, -1.8, -1.391, 2610341.411, 332920971.88596, 1307067, 4.5, 5.5, 4.5, 7, 2, 'bearish', 0.75, 0.09, 'low', 54.49, 'neutral', -0.7195, 'bullish', 0.024, 0.04, 0, 58.4, 'Neutral', 0.477, 0.523, 9.109999999999999, 123.01, -3.57, 'downward', 127.37, 128.24, 127.37, 0.17, 'Hold or wait for clearer signals', 0.54),
    @('20:06:00 21/01/2026', 46043.83750339169, 129.05, 129.77, 124.68, -0.03, -0.023, 2776039.268, 353720271.19413, 1376876, 4.5, 5.5, 4.5, 7, 2, 'bearish', 0.75, 0.23, 'low', 61.65, 'neutral', -0.5184, 'bullish', -0.005, -0.008, 0, 58.9, 'Neutral', 0.474, 0.526, 11.38, 124.75, -3.34, 'downward', 128.88, 129.77, 128.88, 0.17, 'Hold or wait for clearer signals', 0.54),
    @('20:15:46 21/01/2026', 46043.8442848441, 129.09, 129.5, 124.68, -0.09, -0.07000000000000001, 2775722.673, 353672239.4001, 1373155, 4.5, 5.5, 4.5, 7, 2, 'bearish', 0.75, 0.24, 'low', 61.85, 'neutral', -0.5152, 'bullish', -0.005, -0.008, 0, 58.9, 'Neutral', 0.474, 0.526, 11.46, 124.8, -3.33, 'downward', 128.89, 129.77, 128.89, 0.17, 'Hold or wait for clearer signals', 0.54),
    @('20:40:43 21/01/2026', 46043.86161179529, 130.05, 130.33, 124.68, 2.69, 2.112, 2732775.889, 348314476.13076, 1328020, 4.5, 5.5, 4.5, 7, 2, 'neutral', 0.5, 0.32, 'low', 65.95999999999999, 'neutral', -0.2616, 'bullish', -0.005, -0.008999999999999999, 0, 60.5, 'Greed', 0.55, 0.45, 12.46, 126, -3.11, 'downward', 129.7, 131.05, 129.7, 0, 'Hold or wait for clearer signals', 0.55),
    @('22:41:11 21/01/2026', 46043.94527263061, 126.44, 131.18, 124.68, -0.98, -0.769, 2858599.373, 365056668.42764, 1359874, 4.5, 5.5, 4.5, 7, 2, 'bearish', 1, -0.09, 'low', 45.6, 'neutral', -0.4968, 'bullish', -0.006, -0.01, 0, 58.3, 'Neutral', 0.449, 0.551, 20.65, 122.41, -3.18, 'downward', 125.16, 127.64, 125.16, 0.17, 'Hold or wait for clearer signals', 0.55)
)

$dataSheet = $wb.Worksheets.Item("Data")
$startRow = 23
for ($i = 0; $i -lt $dataRows.Count; $i++) {
    $r = $startRow + $i
    $row = $dataRows[$i]
    for ($c = 1; $c -le $row.Count; $c++) {
        $dataSheet.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    # Column B carries the same datetime serial/number-format as the rows above it.
    $dataSheet.Range("B$r").NumberFormat = $dataSheet.Range("B" + ($r - 1)).NumberFormat
}

# ---------------------------------------------------------------------
# 2) Per-topic sheets: each mirrors a subset of the Data columns, written
#    to rows 25-29 (their sheet rows run 2 below the Data sheet's).
# ---------------------------------------------------------------------
$topicSheets = @(
    @{ Name = "Price Analysis";       Cols = @(1, 2, 3, 4, 5, 6, 7, 8) },
    @{ Name = "Technical Analysis";   Cols = @(1, 2, 17, 18, 20, 22) },
    @{ Name = "Fundamental Analysis"; Cols = @(1, 2, 11, 12, 13, 14, 15) },
    @{ Name = "Sentiment Analysis";   Cols = @(1, 2, 24, 25, 26, 27) },
    @{ Name = "Predictions";          Cols = @(1, 2, 3, 32, 33, 29, 30, 35, 36) }
)

foreach ($topic in $topicSheets) {
    $ws = $wb.Worksheets.Item($topic.Name)
    $cols = $topic.Cols
    for ($i = 0; $i -lt $dataRows.Count; $i++) {
        $r = 25 + $i
        $row = $dataRows[$i]
        for ($c = 0; $c -lt $cols.Count; $c++) {
            $ws.Cells.Item($r, $c + 1).Value = $row[$cols[$c] - 1]
        }
        $ws.Range("B$r").NumberFormat = $ws.Range("B" + ($r - 1)).NumberFormat
    }
}

# ---------------------------------------------------------------------
# 3) Dashboard "Latest Analysis" + "Summary Statistics" refresh.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B4").Value = "22:41:11 21/01/2026"
$dash.Range("B5").Value = "$126.44"
$dash.Range("B6").Value = "-0.77%"
$dash.Range("B7").Value = "$122.41"
$dash.Range("B9").Value = "-0.006"
$dash.Range("B10").Value = "0.55"
$dash.Range("B15").Value = 26
$dash.Range("B16").Value = "$127.68"
$dash.Range("B18").Value = "12.35%"
$dash.Range("B19").Value = "37.31"
$dash.Range("B20").Value = "0.022"

# ---------------------------------------------------------------------
# 4) Charts: every series on every chart in the topic sheets currently
#    plots rows 4-24; extend that to 4-29 to reach the new data.
# ---------------------------------------------------------------------
foreach ($topic in $topicSheets) {
    $ws = $wb.Worksheets.Item($topic.Name)
    $chartCount = $ws.ChartObjects().Count
    for ($ci = 1; $ci -le $chartCount; $ci++) {
        $chart = $ws.ChartObjects($ci).Chart
        $seriesCount = $chart.SeriesCollection().Count
        for ($si = 1; $si -le $seriesCount; $si++) {
            $ser = $chart.SeriesCollection($si)
            $ser.Formula = $ser.Formula.Replace("`$4:`$B`$24", "`$4:`$B`$29").Replace("`$4:`$C`$24", "`$4:`$C`$29").Replace("`$4:`$D`$24", "`$4:`$D`$29").Replace("`$4:`$E`$24", "`$4:`$E`$29").Replace("`$4:`$F`$24", "`$4:`$F`$29").Replace("`$4:`$G`$24", "`$4:`$G`$29").Replace("`$4:`$H`$24", "`$4:`$H`$29").Replace("`$4:`$I`$24", "`$4:`$I`$29")
        }
    }
}
